$d = $word.ActiveDocument

# 1) After the "Technology" paragraph, append a new sentence about the
#    Bagging classifier / LightGBM base eliminator.
$d.Content.Find.Execute(
    "is used which is a gradient boosting framework that uses tree-based learning algorithms.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "is used which is a gradient boosting framework that uses tree-based learning algorithms. In addition Bagging classifier is used with LightGBM as a base eliminator. It improves accuracy of model.",
    2) | Out-Null

# Split the newly-added "LightGBM" mention into its own run (mirrors the
# original "Technology : LightGBM ..." sentence, where the product name is
# kept in a dedicated run / proofing span rather than fused into the
# surrounding sentence text).
$lgFind = $d.Content
if ($lgFind.Find.Execute("with LightGBM as")) {
    $lgStart = $lgFind.Start + 5
    $lgEnd = $lgStart + 8
    $lgRange = $d.Range($lgStart, $lgEnd)
    if ($lgRange.Text -eq "LightGBM") {
        $lgRange.Bold = 1
        $lgRange.Bold = 0
    }
}

# 2) Drop "employment_industry" and "employment_occupation" from the list of
#    irrelevant columns, keeping "respondent_id" and "health_insurance" (and
#    their separating ", " run) untouched.
$d.Content.Find.Execute(
    ", employment_industry, employment_occupation, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", ",
    2) | Out-Null

# 3) Tidy up the "evaluate_model" runs/proofing splits (re-typing the same
#    text merges the previously split runs).
$d.Content.Find.Execute(
    "evaluate_model",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "evaluate_model",
    2) | Out-Null
$d.Content.Find.Execute(
    "evaluate_model",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "evaluate_model",
    2) | Out-Null

# 4) Tidy up the "probabilities" run split in item 10).
$d.Content.Find.Execute(
    " and the predicted probabilities for both vaccines is created and saved to a CSV file named results.csv. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and the predicted probabilities for both vaccines is created and saved to a CSV file named results.csv. ",
    2) | Out-Null

# 5) Update the reported ROC AUC result value.
$d.Content.Find.Execute(
    "Result : ROC AUC=0.84822",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Result : ROC AUC=0.851712",
    2) | Out-Null

# 6) Merge the "(-+0.01). This suggests" paragraph with the following
#    "stability in the models' ability..." paragraph into a single
#    paragraph, joined by a space.
$d.Content.Find.Execute(
    "This suggests^pstability",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This suggests stability",
    2) | Out-Null
